$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.136025667190552
$ws.Range("B1").Value = 2.290764093399048
$ws.Range("C1").Value = 11.14107227325439
$ws.Range("D1").Value = 2.139928102493286
$ws.Range("E1").Value = 1.275053381919861
